# Remove the "Max TANF benefit for HH Size:" paragraph (Step 4 / Benefits
# Estimate section). The paragraph - along with its trailing paragraph
# mark - is deleted in its entirety, leaving the preceding blank
# paragraph directly followed by the "Estimated TANF Benefit Amount:"
# paragraph.

$d = $word.ActiveDocument

# Locate the paragraph whose text begins with the target phrase and
# delete its whole Range (text + paragraph mark), which removes the
# paragraph from the document completely.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Max TANF benefit for HH Size:*") {
        $p.Range.Delete()
    }
}
